$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 777; this shifts the existing rows 777-818 down to 778-819
$ws.Rows.Item(777).Insert()

# Populate the newly inserted row 777 with the new data point.
# Column A stores the date as plain text (not an actual Excel date), so we
# prefix with an apostrophe to stop Excel from auto-converting it to a date
# serial number, then clear the resulting formatting so no style gets
# attached (matching the rest of the sheet, which uses the default style).
$ws.Range("A777").Value = "'2026/02/04"
$ws.Range("A777").ClearFormats()
$ws.Range("B777").Value = "水"
$ws.Range("C777").Value = 2
$ws.Range("D777").Value = 201
